# Training Kit deck — "Added training package" pass.
# The author opened the deck, clicked into the Installation & Setup
# bullet list (slide 3 / shape "Content Placeholder 2") and the content
# placeholder got nudged/resized by PowerPoint's autofit while the
# wording itself stayed the same (the GitHub Copilot / RooCode steps).
#
# We reapply the same wording (so spell-check style run boundaries line
# up the same way PowerPoint split them around "licence" / "RooCode")
# and restore the placeholder's resized/auto-fit geometry.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Paragraph "Ensure your GitHub account has a Copilot subscription or
#     enterprise licence" -> split into two runs at the word "licence".
$part = $tr.Characters(16, 68)
$part.Text = "Ensure your GitHub account has a Copilot subscription or enterprise "
$part = $tr.Characters(84, 7)
$part.Text = "licence"

# --- "RooCode" heading run right after the manual line break.
$part = $tr.Characters(240, 7)
$part.Text = "RooCode"

# --- Paragraph "Install RooCode extension from roocode.com" -> split
#     into three runs around the product name "RooCode".
$part = $tr.Characters(248, 8)
$part.Text = "Install "
$part = $tr.Characters(256, 7)
$part.Text = "RooCode"
$part = $tr.Characters(263, 27)
$part.Text = " extension from roocode.com"

# --- Placeholder got resized/repositioned by PowerPoint's "shrink text on
#     overflow" autofit once the content settled; restore that geometry.
$shp.Left = 25.4348031496063
$shp.Top = 126.0
$shp.Width = 658.5651968503937
$shp.Height = 261.78259842519685
